$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- sheet1 ("common"): update upload label and append an "export" row ---
$ws1.Range("B6").Value = "上传"

$ws1.Range("A7").Value = "export"
$ws1.Range("B7").Value = "导出"
$ws1.Range("C7").Value = "Export"
$ws1.Range("D7").Value = "Exporter"

# match formatting of the row above (D column carries a highlighted style)
$ws1.Range("D6").Copy()
$ws1.Range("D7").PasteSpecial(-4122)
$ws1.Rows.Item(7).RowHeight = 14.25

# --- sheet2 ("simulator"): append store_name / store_code translation rows ---
$ws2.Range("A3").Value = "store_name"
$ws2.Range("B3").Value = "商店名称"
$ws2.Range("C3").Value = "Store Name"
$ws2.Range("D3").Value = "Le nom de magasin"

$ws2.Range("A4").Value = "store_code"
$ws2.Range("B4").Value = "商店代码"
$ws2.Range("C4").Value = "Store Code"
$ws2.Range("D4").Value = "Le Code de magasin"

# --- selection / active tab bookkeeping ---
[void]$ws1.Range("C13").Select()
[void]$ws2.Range("B15").Select()
[void]$ws2.Activate()
